# Weekly refresh of Fruta/Hortaliza data: the rows in the sheet have been
# reshuffled (same set of records, new row order). Only the columns that
# vary per-record (D, L, M, N, O, P, R, S) need to move; the columns that
# are constant across every data row (A, B, C, E, F, G, H, I, J, K, Q, T)
# are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column indices for the fields that move with each record.
$cols = @(4, 12, 13, 14, 15, 16, 18, 19)   # D, L, M, N, O, P, R, S

# Mapping of destination data row -> source data row (1-based sheet rows,
# header is row 1, data rows are 2..41). This is the permutation that
# reorders the records.
$map = @{
     2 = 37;  3 = 33;  4 = 28;  5 = 32;  6 = 11;  7 = 21;  8 = 6;   9 = 20
    10 = 41; 11 = 19; 12 = 15; 13 = 23; 14 = 22; 15 = 3;  16 = 30; 17 = 9
    18 = 2;  19 = 38; 20 = 27; 21 = 36; 22 = 17; 23 = 29; 24 = 12; 25 = 40
    26 = 18; 27 = 25; 28 = 34; 29 = 26; 30 = 5;  31 = 16; 32 = 14; 33 = 8
    34 = 24; 35 = 7;  36 = 4;  37 = 31; 38 = 35; 39 = 10; 40 = 13; 41 = 39
}

# Snapshot the original values of the moving columns for every data row
# before any writes happen, so the permutation can be applied safely
# even though source and destination rows overlap. Use .Value2 (not
# .Value) so dates come back as raw serial numbers instead of wrapped
# date variants.
$snapshot = @{}
for ($r = 2; $r -le 41; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowVals
}

# Write the values from the recorded source row into each destination row.
foreach ($destRow in $map.Keys) {
    $srcRow = $map[$destRow]
    $srcVals = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Cells.Item($destRow, $c).Value = $srcVals[$c]
    }
}
